$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text values (prices, percents) are stored as literal text,
# not auto-converted to numbers, matching the original inline-string format.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "62.536.27"
$ws.Range("E2").Value = "  +4.08%  "
$ws.Range("D3").Value = "2.406.99"
$ws.Range("E3").Value = "  +1.22%  "
$ws.Range("E4").Value = "  +0.40%  "
$ws.Range("D5").Value = "573.26"
$ws.Range("E5").Value = "  +2.35%  "
$ws.Range("D6").Value = "145.74"
$ws.Range("E6").Value = "  +5.42%  "
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("E8").Value = "  +2.21%  "
$ws.Range("D9").Value = "2.434.94"
$ws.Range("E9").Value = "  +2.47%  "
$ws.Range("E10").Value = "  +4.64%  "
$ws.Range("E11").Value = "  +0.77%  "
$ws.Range("D12").Value = "5.23"
$ws.Range("E12").Value = "  +2.69%  "
$ws.Range("D13").Value = "0.352"
$ws.Range("E13").Value = "  +4.08%  "
$ws.Range("D14").Value = "27.29"
$ws.Range("E14").Value = "  +6.02%  "
$ws.Range("E15").Value = "  +6.28%  "
$ws.Range("D16").Value = "2.843.60"
$ws.Range("E16").Value = "  +1.76%  "
$ws.Range("D17").Value = "62.597.57"
$ws.Range("E17").Value = "  +4.51%  "
$ws.Range("D18").Value = "2.422.41"
$ws.Range("E18").Value = "  +2.36%  "
$ws.Range("E19").Value = "  -2.31%  "
$ws.Range("D20").Value = "10.95"
$ws.Range("E20").Value = "  +4.15%  "
$ws.Range("D21").Value = "327.37"
$ws.Range("E21").Value = "  +1.71%  "
$ws.Range("E22").Value = "  +1.87%  "
$ws.Range("D23").Value = "2.03"
$ws.Range("E23").Value = "  +11.61%  "
$ws.Range("D24").Value = "0.998"
$ws.Range("E24").Value = "  -0.40%  "
$ws.Range("D25").Value = "65.57"
$ws.Range("E25").Value = "  +2.43%  "
$ws.Range("D26").Value = "627.07"
$ws.Range("E26").Value = "  +12.67%  "
$ws.Range("B27").Value = "Aptos"
$ws.Range("C27").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D27").Value = "8.45"
$ws.Range("E27").Value = "  +4.06%  "
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").Value = "0.0₃0985"
$ws.Range("E28").Value = "  +6.08%  "
$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").Value = "2.529.96"
$ws.Range("E29").Value = "  +1.41%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "8.19"
$ws.Range("E30").Value = "  +2.41%  "
$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").Value = "1.42"
$ws.Range("E31").Value = "  +8.17%  "
$ws.Range("B32").Value = "BabyDogeCoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D32").Value = "0.0₆0443"
$ws.Range("E32").Value = "  +51.62%  "
$ws.Range("D33").Value = "0.137"
$ws.Range("E33").Value = "  +4.51%  "
$ws.Range("D34").Value = "1.84"
$ws.Range("E34").Value = "  +3.17%  "
$ws.Range("E35").Value = "  +3.36%  "
$ws.Range("E36").Value = "  -0.35%  "
$ws.Range("E37").Value = "  +4.45%  "
$ws.Range("D38").Value = "0.374"
$ws.Range("E38").Value = "  +2.02%  "
$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").Value = "151.13"
$ws.Range("E39").Value = "  -0.33%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D40").Value = "5.38"
$ws.Range("E40").Value = "  +7.14%  "
$ws.Range("D41").Value = "18.64"
$ws.Range("E41").Value = "  +2.69%  "
$ws.Range("D42").Value = "2.76"
$ws.Range("E42").Value = "  +14.22%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "1.76"
$ws.Range("E43").Value = "  +6.31%  "
$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "144.66"
$ws.Range("E45").Value = "  +2.70%  "
$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").Value = "3.59"
$ws.Range("E46").Value = "  +2.16%  "
$ws.Range("B47").Value = "WhiteBITCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D47").Value = "14.14"
$ws.Range("E47").Value = "  +20.86%  "
$ws.Range("D48").Value = "20.51"
$ws.Range("E48").Value = "  +7.19%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "0.598"
$ws.Range("E49").Value = "  +1.96%  "
$ws.Range("B50").Value = "Hedera"
$ws.Range("C50").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D50").Value = "0.0515"
$ws.Range("E50").Value = "  +2.93%  "
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").Value = "0.0919"
$ws.Range("E51").Value = "  +2.17%  "
